# Update "想去人数" (F column) figures on the sheets that hold the
# event listing data: "展览" (sheet 1) and "全部类型" (sheet 4).
# Sheets "演出" and "本地生活" only contain header rows, so nothing to
# change there.

$wb = $excel.ActiveWorkbook

# Row -> new F-column value map (old values are shown for reference)
#   Row  Old    New
#    3   104 -> 106
#    4  1586 -> 1592
#    5   604 -> 607
#    6  1090 -> 1091
#    7     5 -> 7
#    8 11360 -> 11366
#   10    89 -> 90
#   11   444 -> 445
#   15 12331 -> 12335
#   16 12990 -> 13002
#   21    82 -> 85

$updates = @{
    3  = 106
    4  = 1592
    5  = 607
    6  = 1091
    7  = 7
    8  = 11366
    10 = 90
    11 = 445
    15 = 12335
    16 = 13002
    21 = 85
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
